$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 7646
$ws.Range("E2").Value = 141
$ws.Range("F2").Value = 141
$ws.Range("G2").Value = 64
$ws.Range("H2").Value = 42
$ws.Range("I2").Value = -51
$ws.Range("J2").Value = 93
$ws.Range("K2").Value = 9848
$ws.Range("L2").Value = 4241
$ws.Range("M2").Value = 5608
$ws.Range("N2").Value = 5191
$ws.Range("O2").Value = 416
$ws.Range("P2").Value = 111
$ws.Range("Q2").Value = 51
$ws.Range("R2").Value = -107
$ws.Range("S2").Value = 132
$ws.Range("T2").Value = 317
$ws.Range("U2").Value = -266
$ws.Range("V2").Value = 2680
$ws.Range("W2").Value = 1.84
$ws.Range("X2").Value = 0.55
$ws.Range("Y2").Value = -0.97
$ws.Range("Z2").Value = 0.43
$ws.Range("AA2").Value = 75.62
$ws.Range("AB2").Value = 4860.3
$ws.Range("AC2").Value = -2059
$ws.Range("AD2").Value = -32.9
$ws.Range("AE2").Value = 279713
$ws.Range("AF2").Value = 0.25
$ws.Range("AG2").Value = 449
$ws.Range("AH2").Value = 0.66
$ws.Range("AI2").Value = -15.8
$ws.Range("AJ2").Value = 2465609

# Row 3
$ws.Range("D3").Value = 7741
$ws.Range("E3").Value = 104
$ws.Range("F3").Value = 104
$ws.Range("G3").Value = -17
$ws.Range("H3").Value = -75
$ws.Range("I3").Value = -148
$ws.Range("J3").Value = 73
$ws.Range("K3").Value = 9787
$ws.Range("L3").Value = 4321
$ws.Range("M3").Value = 5466
$ws.Range("N3").Value = 5023
$ws.Range("O3").Value = 444
$ws.Range("P3").Value = 113
$ws.Range("Q3").Value = 257
$ws.Range("R3").Value = -223
$ws.Range("S3").Value = 195
$ws.Range("T3").Value = 376
$ws.Range("U3").Value = -119
$ws.Range("V3").Value = 2958
$ws.Range("W3").Value = 1.35
$ws.Range("X3").Value = -0.96
$ws.Range("Y3").Value = -2.89
$ws.Range("Z3").Value = -0.76
$ws.Range("AA3").Value = 79.05
$ws.Range("AB3").Value = 4624.38
$ws.Range("AC3").Value = -5984
$ws.Range("AD3").Value = -9.7
$ws.Range("AE3").Value = 274436
$ws.Range("AF3").Value = 0.22
$ws.Range("AG3").Value = 457
$ws.Range("AH3").Value = 0.79
$ws.Range("AI3").Value = -5.48
$ws.Range("AJ3").Value = 2465609

# Row 4
$ws.Range("D4").Value = 8224
$ws.Range("E4").Value = 173
$ws.Range("F4").Value = 173
$ws.Range("G4").Value = 90
$ws.Range("H4").Value = 11
$ws.Range("I4").Value = -38
$ws.Range("J4").Value = 48
$ws.Range("K4").Value = 9869
$ws.Range("L4").Value = 4445
$ws.Range("M4").Value = 5424
$ws.Range("N4").Value = 4956
$ws.Range("O4").Value = 468
$ws.Range("P4").Value = 115
$ws.Range("Q4").Value = 86
$ws.Range("R4").Value = 5
$ws.Range("S4").Value = -87
$ws.Range("T4").Value = 119
$ws.Range("U4").Value = -33
$ws.Range("V4").Value = 2984
$ws.Range("W4").Value = 2.1
$ws.Range("X4").Value = 0.13
$ws.Range("Y4").Value = -0.75
$ws.Range("Z4").Value = 0.11
$ws.Range("AA4").Value = 81.95
$ws.Range("AB4").Value = 4502.21
$ws.Range("AC4").Value = -1523
$ws.Range("AD4").Value = -37.37
$ws.Range("AE4").Value = 270773
$ws.Range("AF4").Value = 0.21
$ws.Range("AG4").Value = 465
$ws.Range("AH4").Value = 0.82
$ws.Range("AI4").Value = -22.08
$ws.Range("AJ4").Value = 2465609

# Row 5
$ws.Range("D5").Value = 8471
$ws.Range("E5").Value = 220
$ws.Range("F5").Value = 220
$ws.Range("G5").Value = -8
$ws.Range("H5").Value = -41
$ws.Range("I5").Value = -63
$ws.Range("J5").Value = 22
$ws.Range("K5").Value = 9311
$ws.Range("L5").Value = 4100
$ws.Range("M5").Value = 5212
$ws.Range("N5").Value = 4854
$ws.Range("O5").Value = 358
$ws.Range("P5").Value = 117
$ws.Range("Q5").Value = 355
$ws.Range("R5").Value = -129
$ws.Range("S5").Value = -404
$ws.Range("T5").Value = 223
$ws.Range("U5").Value = 132
$ws.Range("V5").Value = 2615
$ws.Range("W5").Value = 2.59
$ws.Range("X5").Value = -0.48
$ws.Range("Y5").Value = -1.29
$ws.Range("Z5").Value = -0.43
$ws.Range("AA5").Value = 78.66
$ws.Range("AB5").Value = 4361.32
$ws.Range("AC5").Value = -2558
$ws.Range("AD5").Value = -20.28
$ws.Range("AE5").Value = 265111
$ws.Range("AF5").Value = 0.2
$ws.Range("AG5").Value = 474
$ws.Range("AH5").Value = 0.91
$ws.Range("AI5").Value = -13.48
$ws.Range("AJ5").Value = 2465609

# Row 6
$ws.Range("D6").Value = 9093
$ws.Range("E6").Value = 445
$ws.Range("F6").Value = 445
$ws.Range("G6").Value = 266
$ws.Range("H6").Value = 138
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 9476
$ws.Range("L6").Value = 4350
$ws.Range("M6").Value = 5126
$ws.Range("N6").Value = 4806
$ws.Range("P6").Value = 119
$ws.Range("Q6").Value = 296
$ws.Range("R6").Value = -128
$ws.Range("S6").Value = -128
$ws.Range("T6").Value = 173
$ws.Range("U6").Value = 123
$ws.Range("V6").Value = 2736
$ws.Range("W6").Value = 4.89
$ws.Range("X6").Value = 1.52
$ws.Range("Y6").Value = 0.01
$ws.Range("Z6").Value = 1.47
$ws.Range("AA6").Value = 84.86
$ws.Range("AB6").Value = 4269.26
$ws.Range("AC6").Value = 10
$ws.Range("AD6").Value = 6007.44
$ws.Range("AE6").Value = 265262
$ws.Range("AF6").Value = 0.22
$ws.Range("AG6").Value = 1447
$ws.Range("AH6").Value = 2.47
$ws.Range("AI6").Value = 10778.25
$ws.Range("AJ6").Value = 2465609

# Clear rows 7-9 data cells (D through AI), leaving A/B/C intact
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()